$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view cosmetics -------------------------------------------------
# Best-effort: adjust the workbook window geometry and un-minimize it.
$win = $excel.ActiveWindow
$win.WindowState = -4143   # xlNormal (un-minimize)
$win.Top = 460
$win.Width = 33600
$win.Height = 19060

# Sheet-level zoom for the active sheet view.
$win.Zoom = 150

# --- New column F: "IsRelevant" ---------------------------------------------
# Header cell, formatted like the other header cells (font/fill from C1) but
# with a left+right-only thin border instead of the full box border.
$ws.Range("F1").Value = "IsRelevant"
$ws.Range("C1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Borders.Item(8).LineStyle = -4142
$ws.Range("F1").Borders.Item(9).LineStyle = -4142

# Body cells F2:F15 — same (empty) style as the rest of the data rows (E2's
# style), applied across the whole column range in one paste.
$ws.Range("E2").Copy()
$ws.Range("F2:F15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Move the selection -------------------------------------------------
$ws.Range("C2").Select()
